# CASH REPORT.xlsx edit
# 1) Replace BATCHID value 915726 -> 250080 on the three existing sheets.
# 2) Add four new worksheets (with header rows only) after Batch_Miscellaneous:
#      Batch_Sub_Detail, Batch_Detail_Adjustments,
#      Batch_Header_Optional_Fields, Batch_Detail_Optional_Fields

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) BATCHID value swap: 915726 -> 250080
# ---------------------------------------------------------------------------

$oldBatchId = 915726
$newBatchId = 250080

$batchHeader = $wb.Worksheets.Item("Batch_Header")
for ($r = 2; $r -le 5; $r++) {
    $batchHeader.Cells.Item($r, 2).Value = $newBatchId
}

$batchDetail = $wb.Worksheets.Item("Batch_Detail")
for ($r = 2; $r -le 97; $r++) {
    $batchDetail.Cells.Item($r, 3).Value = $newBatchId
}

$batchMisc = $wb.Worksheets.Item("Batch_Miscellaneous")
for ($r = 2; $r -le 5; $r++) {
    $batchMisc.Cells.Item($r, 1).Value = $newBatchId
}

# ---------------------------------------------------------------------------
# 2) New worksheets (CSV-template based DataFrames), header row only.
# ---------------------------------------------------------------------------

# Use the existing (bold/bordered/centered) header style from Batch_Header!A1
# as the style source so the new header cells share the same style index.
$styleSource = $wb.Worksheets.Item("Batch_Header")

function Add-TemplateSheet {
    param(
        [string]$SheetName,
        [string[]]$Headers
    )

    $lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
    $newWs = $wb.Worksheets.Add($null, $lastSheet)
    $newWs.Name = $SheetName

    $srcRange = $styleSource.Range($styleSource.Cells.Item(1, 1), $styleSource.Cells.Item(1, $Headers.Length))
    $dstRange = $newWs.Range($newWs.Cells.Item(1, 1), $newWs.Cells.Item(1, $Headers.Length))
    $srcRange.Copy($dstRange)

    for ($i = 0; $i -lt $Headers.Length; $i++) {
        $newWs.Cells.Item(1, $i + 1).Value = $Headers[$i]
    }
}

$subDetailHeaders = @(
    "BATCHID","ENTRYNO","DETAILNO","DOCNUMBER","PAYNUMBER","SUBDETNO","DOCTYPE",
    "APPLAMOUNT","DISCOUNT","DATEDOC","IDCUST","PONUMBER","NATCUSTID","ADJREF",
    "ADJDESC","ADJAMOUNT","SWJOB","SWPOSTED","PJCSDAMT","PJCSDDISC","ENTRYTYPE",
    "PROCESSCMD","PJCUNAPAMT","PJCUNAPDSC","AMTWHD1TC","AMTWHD2TC","AMTWHD3TC",
    "AMTWHD4TC","AMTWHD5TC","AMTWHD1HC","AMTWHD2HC","AMTWHD3HC","AMTWHD4HC",
    "AMTWHD5HC","AMTWHD1BC","AMTWHD2BC","AMTWHD3BC","AMTWHD4BC","AMTWHD5BC",
    "AMTWHTOTBC","AMTWHD1DT","AMTWHD2DT","AMTWHD3DT","AMTWHD4DT","AMTWHD5DT",
    "CODETAX1","CODETAX2","CODETAX3","CODETAX4","CODETAX5","AMTWHDTOT",
    "CBBTADRVH","CBBTADVW","AMTREMAIN","APPLYMETH","PNDADJAMT","SELECTFROM"
)

$detailAdjustmentsHeaders = @(
    "BATCHID","ENTRYNO","DETAILNO","DOCNUMBER","PAYNUMBER","TRANSTYPE","SEQNO",
    "DISTCODE","DISTAMOUNT","ACCTID","ACCTIDUF","CONTRACT","PROJECT","CATEGORY",
    "RESOURCE","COSTCLASS","BILLTYPE","ITEMNO","UNITOFMEAS","QUANTITY","COST",
    "DATEBILL","BILLRATE","BILLCUR","FMTCONTNO","DOCLINE","CBBTHDRVH","CBBTHDVW",
    "SWFROMWEB","AMTREMAIN"
)

$headerOptionalFieldsHeaders = @(
    "BATCHID","ENTRYNO","OPTFIELD","VALUE","TYPE","LENGTH","DECIMALS","ALLOWNULL",
    "VALIDATE","SWSET","VALINDEX","VALIFTEXT","VALIFMONEY","VALIFNUM","VALIFLONG",
    "VALIFBOOL","VALIFDATE","VALIFTIME","FDESC","VDESC"
)

$detailOptionalFieldsHeaders = @(
    "BATCHID","ENTRYNO","DETAILNO","OPTFIELD","VALUE","TYPE","LENGTH","DECIMALS",
    "ALLOWNULL","VALIDATE","SWSET","VALINDEX","VALIFTEXT","VALIFMONEY","VALIFNUM",
    "VALIFLONG","VALIFBOOL","VALIFDATE","VALIFTIME","FDESC","VDESC"
)

Add-TemplateSheet -SheetName "Batch_Sub_Detail" -Headers $subDetailHeaders
Add-TemplateSheet -SheetName "Batch_Detail_Adjustments" -Headers $detailAdjustmentsHeaders
Add-TemplateSheet -SheetName "Batch_Header_Optional_Fields" -Headers $headerOptionalFieldsHeaders
Add-TemplateSheet -SheetName "Batch_Detail_Optional_Fields" -Headers $detailOptionalFieldsHeaders
